$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ProductName"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "subCategoryName"
$ws.Range("D1").Value = "Brand"
$ws.Range("E1").Value = "Unit"
$ws.Range("F1").Value = "SKU"
$ws.Range("G1").Value = "minimumqty"
$ws.Range("H1").Value = "Qty"
$ws.Range("I1").Value = "description"
$ws.Range("J1").Value = "tax"
$ws.Range("K1").Value = "discount"
$ws.Range("L1").Value = "price"
$ws.Range("M1").Value = "status"
$ws.Range("N1").Value = "img"

$ws.Range("N1").Select()
